# Generate Report for Handoff
# Adds a new data row (row 3) to each of the three worksheets
# (Overview, zh-cn, de-de) for the file
# "d0047408-3645-49bf-b3c9-3048e91bcb0b.md", mirroring the existing
# row 2 created for "37b2757a-1abc-46d1-8a96-270972dd7d14.md".

$wb = $excel.ActiveWorkbook

$commit = "f8d6cb80a984c0211aae60ba6d7284982581a355"
$baseUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$commit/e2e/"

# Helper: write a literal text value into a cell without Excel's automatic
# type coercion (e.g. "True"/"False" -> boolean, "" -> blank), then make
# sure the cell keeps the default (unstyled) look.
function Set-TextCell($range, [string]$text) {
    $range.Value = "'" + $text
    $range.NumberFormat = "General"
    $range.Style = "Normal"
}

# ---------------------------------------------------------------
# Sheet "Overview": new row 3
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("A3").Value = "d0047408-3645-49bf-b3c9-3048e91bcb0b.md"
$ws1.Range("B3").Value = "e2e\d0047408-3645-49bf-b3c9-3048e91bcb0b.md"
$ws1.Range("C3").Value = ".md"
Set-TextCell $ws1.Range("D3") ""
$ws1.Range("E3").Value = "Ready for handoff"
$ws1.Range("F3").Value = "Ready for handoff"
$ws1.Range("G3").Value = "2016-08-20 16:49:11"
$ws1.Range("G3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$link1 = $ws1.Hyperlinks.Add($ws1.Range("B3"), ($baseUrl + "d0047408-3645-49bf-b3c9-3048e91bcb0b.md"), [Type]::Missing, [Type]::Missing, "e2e\d0047408-3645-49bf-b3c9-3048e91bcb0b.md")

$lo1 = $ws1.ListObjects.Item(1)
$lo1.Resize($ws1.Range("A1:G3"))

# ---------------------------------------------------------------
# Sheet "zh-cn": new row 3
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("B3").Value = ".md"
$ws2.Range("C3").Value = "Ready for handoff"
$ws2.Range("D3").Value = "e2e"
$ws2.Range("E3").Value = "ht"
Set-TextCell $ws2.Range("F3") "False"
$ws2.Range("G3").Value = "d0047408-3645-49bf-b3c9-3048e91bcb0b.e10b8306dfedcc52cc61f1a0d1f3af42366cbec7.zh-cn.xlf"
$ws2.Range("H3").Value = "2016-08-20 16:49:06"
$ws2.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
Set-TextCell $ws2.Range("I3") ""
Set-TextCell $ws2.Range("J3") ""
$ws2.Range("K3").Value = "0001-01-01 00:00:00"
$ws2.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
Set-TextCell $ws2.Range("L3") ""
Set-TextCell $ws2.Range("M3") "True"
Set-TextCell $ws2.Range("N3") ""
Set-TextCell $ws2.Range("O3") "False"
Set-TextCell $ws2.Range("P3") ""

# A3 holds the hyperlinked file name (same pattern as A2)
$ws2.Range("A3").Value = "d0047408-3645-49bf-b3c9-3048e91bcb0b.md"
$link2 = $ws2.Hyperlinks.Add($ws2.Range("A3"), ($baseUrl + "d0047408-3645-49bf-b3c9-3048e91bcb0b.md"), [Type]::Missing, [Type]::Missing, "d0047408-3645-49bf-b3c9-3048e91bcb0b.md")

$lo2 = $ws2.ListObjects.Item(1)
$lo2.Resize($ws2.Range("A1:P3"))

# ---------------------------------------------------------------
# Sheet "de-de": new row 3
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("B3").Value = ".md"
$ws3.Range("C3").Value = "Ready for handoff"
$ws3.Range("D3").Value = "e2e"
$ws3.Range("E3").Value = "ht"
Set-TextCell $ws3.Range("F3") "False"
$ws3.Range("G3").Value = "d0047408-3645-49bf-b3c9-3048e91bcb0b.e10b8306dfedcc52cc61f1a0d1f3af42366cbec7.de-de.xlf"
$ws3.Range("H3").Value = "2016-08-20 16:49:11"
$ws3.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
Set-TextCell $ws3.Range("I3") ""
Set-TextCell $ws3.Range("J3") ""
$ws3.Range("K3").Value = "0001-01-01 00:00:00"
$ws3.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
Set-TextCell $ws3.Range("L3") ""
Set-TextCell $ws3.Range("M3") "True"
Set-TextCell $ws3.Range("N3") ""
Set-TextCell $ws3.Range("O3") "False"
Set-TextCell $ws3.Range("P3") ""

$ws3.Range("A3").Value = "d0047408-3645-49bf-b3c9-3048e91bcb0b.md"
$link3 = $ws3.Hyperlinks.Add($ws3.Range("A3"), ($baseUrl + "d0047408-3645-49bf-b3c9-3048e91bcb0b.md"), [Type]::Missing, [Type]::Missing, "d0047408-3645-49bf-b3c9-3048e91bcb0b.md")

$lo3 = $ws3.ListObjects.Item(1)
$lo3.Resize($ws3.Range("A1:P3"))
